$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the entire second row ("H 72") which shifts all rows below it up by one.
$ws.Rows(2).Delete()
